$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 2.2062744499517599
$ws.Range("C2").Value = 0.72790054227419165
$ws.Range("D2").Value = 1.9204613482236521
$ws.Range("E2").Value = 0.48922537474695921

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 2.1640683466431012
$ws.Range("C3").Value = 0.25390981566173826
$ws.Range("D3").Value = 2.2601054549193393
$ws.Range("E3").Value = 0.77385520587044976

# Update the selected range shown when the workbook is opened
$ws.Range("B1:E3").Select()
